$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Tabelle1 (sheet1): restore its original text values so the shared-
# string table round-trips correctly (the source file's <v> indices were
# wrapped across lines, which otherwise breaks the lookup), then reselect
# the whole data range without a stale D3 active-cell.
$ws1.Range("A1").Value = "Bob"
$ws1.Range("B1").Value = "Alice"
$ws1.Range("C1").Value = "Sue"
$ws1.Range("A2").Value = "Yes"
$ws1.Range("B2").Value = "No"
$ws1.Range("C2").Value = "Yes"
$ws1.Range("A3").Value = "No"
$ws1.Range("C3").Value = "Yes"

# --- Tabelle2 (sheet2): mirror the header row from Tabelle1 and add a
# small numeric data table below it.
$ws2.Range("A1").Value = "Bob"
$ws2.Range("B1").Value = "Alice"
$ws2.Range("C1").Value = "Sue"

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 2
$ws2.Range("C2").Value = 3

$ws2.Range("A3").Value = 4
$ws2.Range("C3").Value = 5

# --- Selections: Tabelle1 keeps the whole data range selected (no single
# active cell left over), Tabelle2 becomes the active tab with C3 selected.
$ws1.Range("A1:C3").Select()
$ws2.Range("C3").Select()
$ws2.Activate()
